# Generate Report for Handback
# Fills in the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) for the 6cbd272b... file row (row 3)
# on both the "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-22 05:24:53"
$wsZhCn.Range("H3").Value = "2016-03-22 05:25:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-22 05:25:00"
$wsDeDe.Range("H3").Value = "2016-03-22 05:25:48"
